$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 22.79335135855249
$ws.Range("C2").Value = 15.43741065065874
$ws.Range("D2").Value = 7.302821254366672
$ws.Range("F2").Value = 47.89074097796659
$ws.Range("G2").Value = 3.707206416582363
$ws.Range("J2").Value = 11.5082108940864
$ws.Range("N2").Value = 19.39117377730971

# Row 3
$ws.Range("B3").Value = 22.26613143582738
$ws.Range("C3").Value = 14.94054842355924
$ws.Range("D3").Value = 7.303085630920379
$ws.Range("F3").Value = 47.55594301086751
$ws.Range("G3").Value = 3.712104256277519
$ws.Range("J3").Value = 11.50000191463146
$ws.Range("N3").Value = 19.4586858559522

# Row 4
$ws.Range("B4").Value = 21.94384236518762
$ws.Range("C4").Value = 14.63245411624372
$ws.Range("D4").Value = 7.304327700218745
$ws.Range("F4").Value = 47.36439486007571
$ws.Range("G4").Value = 3.715262060286009
$ws.Range("J4").Value = 11.49757132748295
$ws.Range("N4").Value = 19.50217105464794

# Row 5
$ws.Range("B5").Value = 21.81307377836738
$ws.Range("C5").Value = 14.50636413973217
$ws.Range("D5").Value = 7.305106359894628
$ws.Range("F5").Value = 47.28991401092873
$ws.Range("G5").Value = 3.716586908473576
$ws.Range("J5").Value = 11.49723641367262
$ws.Range("N5").Value = 19.52040294775552

# Row 6
$ws.Range("B6").Value = 21.79140015253231
$ws.Range("C6").Value = 14.48540114461753
$ws.Range("D6").Value = 7.305252142055642
$ws.Range("F6").Value = 47.27776389869508
$ws.Range("G6").Value = 3.716809199652118
$ws.Range("J6").Value = 11.49722036858013
$ws.Range("N6").Value = 19.52346122962357

# Row 7
$ws.Range("B7").Value = 21.94207620104571
$ws.Range("C7").Value = 14.63075550409249
$ws.Range("D7").Value = 7.304337096859642
$ws.Range("F7").Value = 47.36337584123745
$ws.Range("G7").Value = 3.715279773493867
$ws.Range("J7").Value = 11.49756415751487
$ws.Range("N7").Value = 19.50241486565969

# Row 8
$ws.Range("B8").Value = 22.61139344005727
$ws.Range("C8").Value = 15.26684669264481
$ws.Range("D8").Value = 7.302688930398167
$ws.Range("F8").Value = 47.772425911616
$ws.Range("G8").Value = 3.708864056063663
$ws.Range("J8").Value = 11.50483849308243
$ws.Range("N8").Value = 19.41403023060413

# Row 9
$ws.Range("B9").Value = 23.92636260938787
$ws.Range("C9").Value = 16.48135500520887
$ws.Range("D9").Value = 7.307977985852292
$ws.Range("F9").Value = 48.68322354153671
$ws.Range("G9").Value = 3.697469306554482
$ws.Range("J9").Value = 11.5398251266268
$ws.Range("N9").Value = 19.25682382801576

# Row 10
$ws.Range("B10").Value = 24.8826015631336
$ws.Range("C10").Value = 17.34281298891355
$ws.Range("D10").Value = 7.316993132530205
$ws.Range("F10").Value = 49.41472953102308
$ws.Range("G10").Value = 3.689809964452722
$ws.Range("J10").Value = 11.57815932727834
$ws.Range("N10").Value = 19.15113164699919

# Row 11
$ws.Range("B11").Value = 25.31326028547427
$ws.Range("C11").Value = 17.72605855468193
$ws.Range("D11").Value = 7.322192597681701
$ws.Range("F11").Value = 49.76010097394045
$ws.Range("G11").Value = 3.686477871755806
$ws.Range("J11").Value = 11.59833162983435
$ws.Range("N11").Value = 19.10517625476613

# Row 12
$ws.Range("B12").Value = 25.47554575280967
$ws.Range("C12").Value = 17.86980052928824
$ws.Range("D12").Value = 7.324317985295744
$ws.Range("F12").Value = 49.89261000057061
$ws.Range("G12").Value = 3.685237795649531
$ws.Range("J12").Value = 11.60636190624724
$ws.Range("N12").Value = 19.0880796350273

# Row 13
$ws.Range("B13").Value = 25.44063266764669
$ws.Range("C13").Value = 17.83890682528696
$ws.Range("D13").Value = 7.323853310902488
$ws.Range("F13").Value = 49.86399656592683
$ws.Range("G13").Value = 3.685503905122695
$ws.Range("J13").Value = 11.60461506320982
$ws.Range("N13").Value = 19.09174809872253

# Row 14
$ws.Range("B14").Value = 25.32662849567729
$ws.Range("C14").Value = 17.73791281097881
$ws.Range("D14").Value = 7.322364327672679
$ws.Range("F14").Value = 49.77096853206643
$ws.Range("G14").Value = 3.6863754157107
$ws.Range("J14").Value = 11.5989844510707
$ws.Range("N14").Value = 19.10376357782786

# Row 15
$ws.Range("B15").Value = 25.25668913897703
$ws.Range("C15").Value = 17.67586673977877
$ws.Range("D15").Value = 7.321472612937913
$ws.Range("F15").Value = 49.71420805057149
$ws.Range("G15").Value = 3.686912064056285
$ws.Range("J15").Value = 11.59558646184138
$ws.Range("N15").Value = 19.11116322205934

# Row 16
$ws.Range("B16").Value = 24.85435474331362
$ws.Range("C16").Value = 17.31758063622737
$ws.Range("D16").Value = 7.316675312987707
$ws.Range("F16").Value = 49.39240489305637
$ws.Range("G16").Value = 3.690030772424894
$ws.Range("J16").Value = 11.57689590003444
$ws.Range("N16").Value = 19.15417770535318

# Row 17
$ws.Range("B17").Value = 24.60629838678854
$ws.Range("C17").Value = 17.09546515800086
$ws.Range("D17").Value = 7.31401262514343
$ws.Range("F17").Value = 49.19815612929432
$ws.Range("G17").Value = 3.691982856365015
$ws.Range("J17").Value = 11.56612901072105
$ws.Range("N17").Value = 19.18110995214037

# Row 18
$ws.Range("B18").Value = 24.46322404359763
$ws.Range("C18").Value = 16.96690488350215
$ws.Range("D18").Value = 7.312584601334329
$ws.Range("F18").Value = 49.08762144768733
$ws.Range("G18").Value = 3.693119977433473
$ws.Range("J18").Value = 11.56019363412917
$ws.Range("N18").Value = 19.19680062270313

# Row 19
$ws.Range("B19").Value = 24.41471844344965
$ws.Range("C19").Value = 16.92324287935597
$ws.Range("D19").Value = 7.312118914491658
$ws.Range("F19").Value = 49.05040366881776
$ws.Range("G19").Value = 3.693507454091785
$ws.Range("J19").Value = 11.5582282791209
$ws.Range("N19").Value = 19.20214755162638

# Row 20
$ws.Range("B20").Value = 24.63274697385637
$ws.Range("C20").Value = 17.11919416287108
$ws.Range("D20").Value = 7.314285373162456
$ws.Range("F20").Value = 49.21871149461483
$ws.Range("G20").Value = 3.691773571399584
$ws.Range("J20").Value = 11.56724852972466
$ws.Range("N20").Value = 19.17822227278223

# Row 21
$ws.Range("B21").Value = 25.36013718748072
$ws.Range("C21").Value = 17.76761584712994
$ws.Range("D21").Value = 7.322797444058359
$ws.Range("F21").Value = 49.79824704480342
$ws.Range("G21").Value = 3.686118843925656
$ws.Range("J21").Value = 11.6006276895529
$ws.Range("N21").Value = 19.10022604223893

# Row 22
$ws.Range("B22").Value = 25.83081990639779
$ws.Range("C22").Value = 18.18326961477323
$ws.Range("D22").Value = 7.329271812029324
$ws.Range("F22").Value = 50.18701440591551
$ws.Range("G22").Value = 3.682549646515053
$ws.Range("J22").Value = 11.62472357280799
$ws.Range("N22").Value = 19.05103291813736

# Row 23
$ws.Range("B23").Value = 25.58009218175035
$ws.Range("C23").Value = 17.96221430318107
$ws.Range("D23").Value = 7.325733454943779
$ws.Range("F23").Value = 49.97863639847314
$ws.Range("G23").Value = 3.684443075792201
$ws.Range("J23").Value = 11.6116551257427
$ws.Range("N23").Value = 19.0771250918479

# Row 24
$ws.Range("B24").Value = 24.62079099799778
$ws.Range("C24").Value = 17.1084689572042
$ws.Range("D24").Value = 7.314161743500255
$ws.Range("F24").Value = 49.2094148517189
$ws.Range("G24").Value = 3.691868142889991
$ws.Range("J24").Value = 11.56674160187903
$ws.Range("N24").Value = 19.17952714798969

# Row 25
$ws.Range("B25").Value = 23.57157433956777
$ws.Range("C25").Value = 16.15749578873322
$ws.Range("D25").Value = 7.305641338623026
$ws.Range("F25").Value = 48.42557523445979
$ws.Range("G25").Value = 3.700425998226585
$ws.Range("J25").Value = 11.52814051816345
$ws.Range("N25").Value = 19.29762898982517

Write-Host "Updated loading_percent values for 380 kV case"